# Generate Report for Handback
# Refreshes the handoff/handback timestamps for the second tracked file
# (507136a8-379e-44e1-9498-351a84a40d30) after a new de-de xliff hand-back
# round-trip, and rolls the "Latest HO Xliff Generate Date" shown on the
# Overview sheet forward to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: new Correspond Handoff Datetime / Correspond Handback DateTime
# for row 3 (507136a8-379e-44e1-9498-351a84a40d30)
$zhcn.Range("H3").Value = "2016-08-23 20:50:44"
$zhcn.Range("K3").Value = "2016-08-23 20:51:05"

# de-de: new Correspond Handoff Datetime / Correspond Handback DateTime
# for row 3 (507136a8-379e-44e1-9498-351a84a40d30)
$dede.Range("H3").Value = "2016-08-23 20:50:49"
$dede.Range("K3").Value = "2016-08-23 20:51:17"

# Overview: Latest HO Xliff Generate Date for the same file rolls forward
$overview.Range("G3").Value = "2016-08-23 20:50:49"
